$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (both columns end up the same width)
$ws.Columns.Item(1).ColumnWidth = 13.625
$ws.Columns.Item(2).ColumnWidth = 13.625

# Update cell values
$ws.Range("A1").Value = -0.02723143472452776
$ws.Range("B1").Value = 0.027231433527292199

$ws.Range("A2").Value = -0.078096088385835302
$ws.Range("B2").Value = 0.078096087170080097

$ws.Range("A3").Value = 0.023733908569864487
$ws.Range("B3").Value = -0.023733909885533144

$ws.Range("A4").Value = 0.044714723019454124
$ws.Range("B4").Value = -0.044714724272755288
